$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ballot #4's choices to B, C, D (was A, B, C) to introduce a tie-break scenario
$ws.Range("B5").Value = "B"
$ws.Range("C5").Value = "C"
$ws.Range("D5").Value = "D"

# Update ballot #7's choices to C, D, undervote (was B, C, D)
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "D"
$ws.Range("D8").Value = "undervote"

# Remove the last ballot row (ballot #10), shrinking the data set by one row
$ws.Rows(11).Delete()

# Update the active selection to reflect where the cursor ended up
$ws.Range("C11").Select()
